# Aggiornamento dati fino a 1/09/2021
# Adds rows 358-366 (dates 44432-44440) to the sheet, extending the
# existing data table with new daily figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(44432, 0, 10, 65.98046978094484),
    @(44433, 0, 10, 65.98046978094484),
    @(44434, 0, 10, 65.98046978094484),
    @(44435, 3, 8, 52.78437582475588),
    @(44436, 1, 9, 59.38242280285036),
    @(44437, 10, 14, 92.37265769332278),
    @(44438, 0, 14, 92.37265769332278),
    @(44439, 10, 24, 158.3531274742676),
    @(44440, 0, 24, 158.3531274742676)
)

$startRow = 358

# Copy the formatting (styles) of the last existing row down onto the new
# rows before filling in the values, so the new cells match the look of
# the existing table (e.g. the bordered/bold date style in column A).
$lastRow = $startRow - 1
$lastRowRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, 4))
$lastRowRange.Copy()

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $startRow + $i
    $destRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 4))
    $destRange.PasteSpecial(-4122)  # xlPasteFormats
}

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
